# Update betting odds values for rows 2, 3, 5, 8 and 14 to match the
# 2024-10-16 FlashScore refresh (Atualizando o arquivo XLSX).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.7
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.4
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 13
$ws.Range("AB2").Value = 34
$ws.Range("AD2").Value = 6.5
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 26
$ws.Range("AJ2").Value = 19
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 9.5
$ws.Range("AS2").Value = 201
$ws.Range("AW2").Value = 6.5
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 151
$ws.Range("BB2").Value = 401
$ws.Range("G3").Value = 1.8
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.5
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 9
$ws.Range("Q3").Value = 2.1
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 1.7
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 17
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AH3").Value = 10
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 15
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 10
$ws.Range("AQ3").Value = 34
$ws.Range("AX3").Value = 26
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 101
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 301
$ws.Range("G5").Value = 1.55
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 6.5
$ws.Range("L5").Value = 7.5
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 1.44
$ws.Range("W5").Value = 4.75
$ws.Range("X5").Value = 6
$ws.Range("Y5").Value = 9.5
$ws.Range("Z5").Value = 10
$ws.Range("AC5").Value = 6.5
$ws.Range("AD5").Value = 8
$ws.Range("AH5").Value = 12
$ws.Range("AQ5").Value = 29
$ws.Range("AW5").Value = 8
$ws.Range("AZ5").Value = 201
$ws.Range("G8").Value = 2.45
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 2.42
$ws.Range("J8").Value = 2.9
$ws.Range("K8").Value = 2.32
$ws.Range("L8").Value = 2.9
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 13.4
$ws.Range("Q8").Value = 1.52
$ws.Range("R8").Value = 2.2
$ws.Range("U8").Value = 1.47
$ws.Range("V8").Value = 2.32
$ws.Range("W8").Value = 12
$ws.Range("Z8").Value = 27
$ws.Range("AA8").Value = 17.5
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 15.5
$ws.Range("AD8").Value = 7.6
$ws.Range("AG8").Value = 200
$ws.Range("AH8").Value = 11.75
$ws.Range("AI8").Value = 14.5
$ws.Range("AK8").Value = 27
$ws.Range("AL8").Value = 17.5
$ws.Range("AN8").Value = 4.7
$ws.Range("AO8").Value = 12
$ws.Range("AP8").Value = 16.5
$ws.Range("AQ8").Value = 45
$ws.Range("AR8").Value = 60
$ws.Range("AS8").Value = 150
$ws.Range("AT8").Value = 3
$ws.Range("AV8").Value = 40
$ws.Range("AW8").Value = 4.7
$ws.Range("AX8").Value = 12
$ws.Range("BA8").Value = 65
$ws.Range("G14").Value = 2.27
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.85
$ws.Range("J14").Value = 2.9
$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 3.45
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 6.8
$ws.Range("O14").Value = 1.35
$ws.Range("P14").Value = 2.92
$ws.Range("Q14").Value = 2.05
$ws.Range("R14").Value = 1.72
$ws.Range("T14").Value = 2.67
$ws.Range("U14").Value = 1.83
$ws.Range("V14").Value = 1.87
$ws.Range("W14").Value = 7.4
$ws.Range("X14").Value = 10.75
$ws.Range("Y14").Value = 9.25
$ws.Range("Z14").Value = 23
$ws.Range("AA14").Value = 19.5
$ws.Range("AC14").Value = 6.8
$ws.Range("AD14").Value = 6.3
$ws.Range("AE14").Value = 15
$ws.Range("AG14").Value = 700
$ws.Range("AH14").Value = 8.25
$ws.Range("AI14").Value = 14
$ws.Range("AJ14").Value = 10.75
$ws.Range("AK14").Value = 35
$ws.Range("AL14").Value = 26
$ws.Range("AM14").Value = 37
$ws.Range("AN14").Value = 4.15
$ws.Range("AO14").Value = 12
$ws.Range("AP14").Value = 21
$ws.Range("AQ14").Value = 50
$ws.Range("AR14").Value = 90
$ws.Range("AT14").Value = 2.67
$ws.Range("AU14").Value = 7.3
$ws.Range("AV14").Value = 70
$ws.Range("AW14").Value = 4.75
$ws.Range("AX14").Value = 16
$ws.Range("AY14").Value = 25
$ws.Range("AZ14").Value = 75
$ws.Range("BA14").Value = 120
$ws.Range("BB14").Value = 350
